# Test for multiple sheets: add "Test Sheet 2" after "Test Sheet 1",
# populate it with data, and make it the active sheet/tab.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.ActiveSheet

# Add the new sheet right after "Test Sheet 1" and name it.
$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "Test Sheet 2"

# Populate the new sheet with data (matches the target dimension C6:D9).
$ws2.Range("C6").Value = "Random Key"
$ws2.Range("D6").Value = "Random Value"
$ws2.Range("C9").Value = "Random Key 2"
$ws2.Range("D9").Value = 234

# Column widths to roughly match the source template (C ~15.79, D ~15.27 chars).
$ws2.Columns.Item(3).ColumnWidth = 15.0
$ws2.Columns.Item(4).ColumnWidth = 14.5

# Match page setup / header-footer with the rest of the workbook.
$ws2.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ws2.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12Page &P'
$ws2.PageSetup.LeftMargin = 56.7
$ws2.PageSetup.RightMargin = 56.7
$ws2.PageSetup.TopMargin = 75.8
$ws2.PageSetup.BottomMargin = 75.8
$ws2.PageSetup.HeaderMargin = 56.7
$ws2.PageSetup.FooterMargin = 56.7

# Select/activate cell C10 (below the data) and make this the active sheet/tab.
$ws2.Activate()
$ws2.Range("C10").Select() | Out-Null
